$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("A2").Value = 12497
$ws.Range("B2").Value = "Kamilly da Paz"
$ws.Range("C2").Value = "Recursos Humanos"
$ws.Range("D2").Value = "Doenca"
$ws.Range("G2").Value = 2961.18

# Row 3
$ws.Range("A3").Value = 2424
$ws.Range("B3").Value = "João Miguel da Cruz"
$ws.Range("C3").Value = "Atendimento ao Cliente"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45090
$ws.Range("G3").Value = 4964.4

# Row 4
$ws.Range("A4").Value = 50837
$ws.Range("B4").Value = "Sabrina da Conceição"
$ws.Range("C4").Value = "P&D"
$ws.Range("D4").Value = "Viagem de negocios"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 45093
$ws.Range("G4").Value = 5782.89

# Row 5
$ws.Range("A5").Value = 69987
$ws.Range("B5").Value = "Yuri Caldeira"
$ws.Range("C5").Value = "Recursos Humanos"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 45089
$ws.Range("G5").Value = 4441.92

# Row 6
$ws.Range("A6").Value = 70122
$ws.Range("B6").Value = "Gustavo Henrique Guerra"
$ws.Range("C6").Value = "Financeiro"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 8
$ws.Range("G6").Value = 6003.14

# Row 7
$ws.Range("A7").Value = 70687
$ws.Range("B7").Value = "Hadassa Aparecida"
$ws.Range("C7").Value = "Vendas"
$ws.Range("D7").Value = "Viagem de negocios"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45082
$ws.Range("G7").Value = 9820.889999999999

# Row 8
$ws.Range("A8").Value = 65893
$ws.Range("B8").Value = "Fernando Pacheco"
$ws.Range("C8").Value = "Marketing"
$ws.Range("D8").Value = "Viagem de negocios"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45085
$ws.Range("G8").Value = 5551.69

# Row 9
$ws.Range("A9").Value = 68289
$ws.Range("B9").Value = "Dra. Manuella Lima"
$ws.Range("C9").Value = "TI"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45098
$ws.Range("G9").Value = 8131.13

# Row 10
$ws.Range("A10").Value = 6605
$ws.Range("B10").Value = "Luísa Silveira"
$ws.Range("C10").Value = "Recursos Humanos"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("F10").Value = 45100
$ws.Range("G10").Value = 6411.61

# Row 11
$ws.Range("A11").Value = 89521
$ws.Range("B11").Value = "Asafe Leão"
$ws.Range("C11").Value = "Engenharia"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45088
$ws.Range("G11").Value = 8645.450000000001

$wb.Save()
